$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 29 (shifts Sensors SDA..FTDI_RTS and the COUNTA formula row down by one)
$ws.Rows.Item(29).Insert()

# Fill in the new pin mapping row
$ws.Range("A29").Value = "NRF_RESET"
$ws.Range("B29").Value = "PA17"

# Restore the selection to match the committed workbook state
$ws.Range("B29").Select()
